$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 16,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 28.48226033333333
$arr[0,3] = 85.44678099999999
$arr[0,4] = 0.2101651977164657
$arr[0,5] = 0.2101651977164658
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 6.360972666666666
$arr[0,9] = 19.082918
$arr[0,10] = 0.03115862208643261
$arr[0,11] = 0.03115862208643262
$arr[0,12] = 181.1748794652175
$arr[0,13] = 1630.573915186958
$arr[0,14] = 0.006548457971367747
$arr[0,15] = 0.006548457971367749
$arr[1,0] = 3
$arr[1,1] = 1
$arr[1,2] = 28.48226033333333
$arr[1,3] = 85.44678099999999
$arr[1,4] = 0.2101651977164657
$arr[1,5] = 0.2101651977164658
$arr[1,6] = 3
$arr[1,7] = 1
$arr[1,8] = 107.3681206666667
$arr[1,9] = 322.104362
$arr[1,10] = 0.5259325690101214
$arr[1,11] = 0.5259325690101214
$arr[1,12] = 3058.086764328747
$arr[1,13] = 27522.78087895872
$arr[1,14] = 0.1105327223515409
$arr[1,15] = 0.1105327223515409
$arr[2,0] = 3
$arr[2,1] = 1
$arr[2,2] = 28.48226033333333
$arr[2,3] = 85.44678099999999
$arr[2,4] = 0.2101651977164657
$arr[2,5] = 0.2101651977164658
$arr[2,6] = 3
$arr[2,7] = 1
$arr[2,8] = 59.10257333333334
$arr[2,9] = 177.30772
$arr[2,10] = 0.2895083571855736
$arr[2,11] = 0.2895083571855736
$arr[2,12] = 1683.374880049924
$arr[2,13] = 15150.37392044932
$arr[2,14] = 0.06084458112847525
$arr[2,15] = 0.06084458112847526
$arr[3,0] = 3
$arr[3,1] = 1
$arr[3,2] = 28.48226033333333
$arr[3,3] = 85.44678099999999
$arr[3,4] = 0.2101651977164657
$arr[3,5] = 0.2101651977164658
$arr[3,6] = 3
$arr[3,7] = 1
$arr[3,8] = 31.31640666666667
$arr[3,9] = 93.94922000000001
$arr[3,10] = 0.1534004517178723
$arr[3,11] = 0.1534004517178724
$arr[3,12] = 891.9620473845355
$arr[3,13] = 8027.65842646082
$arr[3,14] = 0.0322394362650818
$arr[3,15] = 0.03223943626508181
$arr[4,0] = 3
$arr[4,1] = 1
$arr[4,2] = 83.45109033333334
$arr[4,3] = 250.353271
$arr[4,4] = 0.6157697701763504
$arr[4,5] = 0.6157697701763504
$arr[4,6] = 3
$arr[4,7] = 1
$arr[4,8] = 6.360972666666666
$arr[4,9] = 19.082918
$arr[4,10] = 0.03115862208643261
$arr[4,11] = 0.03115862208643262
$arr[4,12] = 530.8301046138643
$arr[4,13] = 4777.470941524778
$arr[4,14] = 0.01918653756117437
$arr[4,15] = 0.01918653756117437
$arr[5,0] = 3
$arr[5,1] = 1
$arr[5,2] = 83.45109033333334
$arr[5,3] = 250.353271
$arr[5,4] = 0.6157697701763504
$arr[5,5] = 0.6157697701763504
$arr[5,6] = 3
$arr[5,7] = 1
$arr[5,8] = 107.3681206666667
$arr[5,9] = 322.104362
$arr[5,10] = 0.5259325690101214
$arr[5,11] = 0.5259325690101214
$arr[5,12] = 8959.986736674236
$arr[5,13] = 80639.88063006812
$arr[5,14] = 0.32385337714762
$arr[5,15] = 0.32385337714762
$arr[6,0] = 3
$arr[6,1] = 1
$arr[6,2] = 83.45109033333334
$arr[6,3] = 250.353271
$arr[6,4] = 0.6157697701763504
$arr[6,5] = 0.6157697701763504
$arr[6,6] = 3
$arr[6,7] = 1
$arr[6,8] = 59.10257333333334
$arr[6,9] = 177.30772
$arr[6,10] = 0.2895083571855736
$arr[6,11] = 0.2895083571855736
$arr[6,12] = 4932.174186172459
$arr[6,13] = 44389.56767555213
$arr[6,14] = 0.1782704945682934
$arr[6,15] = 0.1782704945682934
$arr[7,0] = 3
$arr[7,1] = 1
$arr[7,2] = 83.45109033333334
$arr[7,3] = 250.353271
$arr[7,4] = 0.6157697701763504
$arr[7,5] = 0.6157697701763504
$arr[7,6] = 3
$arr[7,7] = 1
$arr[7,8] = 31.31640666666667
$arr[7,9] = 93.94922000000001
$arr[7,10] = 0.1534004517178723
$arr[7,11] = 0.1534004517178724
$arr[7,12] = 2613.388281655403
$arr[7,13] = 23520.49453489862
$arr[7,14] = 0.09445936089926259
$arr[7,15] = 0.0944593608992626
$arr[8,0] = 3
$arr[8,1] = 1
$arr[8,2] = 14.78130366666667
$arr[8,3] = 44.343911
$arr[8,4] = 0.1090684366779874
$arr[8,5] = 0.1090684366779875
$arr[8,6] = 3
$arr[8,7] = 1
$arr[8,8] = 6.360972666666666
$arr[8,9] = 19.082918
$arr[8,10] = 0.03115862208643261
$arr[8,11] = 0.03115862208643262
$arr[8,12] = 94.02346860136643
$arr[8,13] = 846.211217412298
$arr[8,14] = 0.003398422200007417
$arr[8,15] = 0.003398422200007418
$arr[9,0] = 3
$arr[9,1] = 1
$arr[9,2] = 14.78130366666667
$arr[9,3] = 44.343911
$arr[9,4] = 0.1090684366779874
$arr[9,5] = 0.1090684366779875
$arr[9,6] = 3
$arr[9,7] = 1
$arr[9,8] = 107.3681206666667
$arr[9,9] = 322.104362
$arr[9,10] = 0.5259325690101214
$arr[9,11] = 0.5259325690101214
$arr[9,12] = 1587.040795693309
$arr[9,13] = 14283.36716123978
$arr[9,14] = 0.0573626430999717
$arr[9,15] = 0.0573626430999717
$arr[10,0] = 3
$arr[10,1] = 1
$arr[10,2] = 14.78130366666667
$arr[10,3] = 44.343911
$arr[10,4] = 0.1090684366779874
$arr[10,5] = 0.1090684366779875
$arr[10,6] = 3
$arr[10,7] = 1
$arr[10,8] = 59.10257333333334
$arr[10,9] = 177.30772
$arr[10,10] = 0.2895083571855736
$arr[10,11] = 0.2895083571855736
$arr[10,12] = 873.6130839214356
$arr[10,13] = 7862.517755292921
$arr[10,14] = 0.0315762239234429
$arr[10,15] = 0.03157622392344291
$arr[11,0] = 3
$arr[11,1] = 1
$arr[11,2] = 14.78130366666667
$arr[11,3] = 44.343911
$arr[11,4] = 0.1090684366779874
$arr[11,5] = 0.1090684366779875
$arr[11,6] = 3
$arr[11,7] = 1
$arr[11,8] = 31.31640666666667
$arr[11,9] = 93.94922000000001
$arr[11,10] = 0.1534004517178723
$arr[11,11] = 0.1534004517178724
$arr[11,12] = 462.8973166888244
$arr[11,13] = 4166.075850199421
$arr[11,14] = 0.01673114745456543
$arr[11,15] = 0.01673114745456544
$arr[12,0] = 3
$arr[12,1] = 1
$arr[12,2] = 8.808546666666667
$arr[12,3] = 26.42564
$arr[12,4] = 0.0649965954291964
$arr[12,5] = 0.06499659542919642
$arr[12,6] = 3
$arr[12,7] = 1
$arr[12,8] = 6.360972666666666
$arr[12,9] = 19.082918
$arr[12,10] = 0.03115862208643261
$arr[12,11] = 0.03115862208643262
$arr[12,12] = 56.03092457972444
$arr[12,13] = 504.27832121752
$arr[12,14] = 0.002025204353883084
$arr[12,15] = 0.002025204353883085
$arr[13,0] = 3
$arr[13,1] = 1
$arr[13,2] = 8.808546666666667
$arr[13,3] = 26.42564
$arr[13,4] = 0.0649965954291964
$arr[13,5] = 0.06499659542919642
$arr[13,6] = 3
$arr[13,7] = 1
$arr[13,8] = 107.3681206666667
$arr[13,9] = 322.104362
$arr[13,10] = 0.5259325690101214
$arr[13,11] = 0.5259325690101214
$arr[13,12] = 945.7571014046313
$arr[13,13] = 8511.813912641681
$arr[13,14] = 0.03418382641098878
$arr[13,15] = 0.03418382641098879
$arr[14,0] = 3
$arr[14,1] = 1
$arr[14,2] = 8.808546666666667
$arr[14,3] = 26.42564
$arr[14,4] = 0.0649965954291964
$arr[14,5] = 0.06499659542919642
$arr[14,6] = 3
$arr[14,7] = 1
$arr[14,8] = 59.10257333333334
$arr[14,9] = 177.30772
$arr[14,10] = 0.2895083571855736
$arr[14,11] = 0.2895083571855736
$arr[14,12] = 520.6077753267556
$arr[14,13] = 4685.469977940801
$arr[14,14] = 0.01881705756536201
$arr[14,15] = 0.01881705756536201
$arr[15,0] = 3
$arr[15,1] = 1
$arr[15,2] = 8.808546666666667
$arr[15,3] = 26.42564
$arr[15,4] = 0.0649965954291964
$arr[15,5] = 0.06499659542919642
$arr[15,6] = 3
$arr[15,7] = 1
$arr[15,8] = 31.31640666666667
$arr[15,9] = 93.94922000000001
$arr[15,10] = 0.1534004517178723
$arr[15,11] = 0.1534004517178724
$arr[15,12] = 275.8520295556444
$arr[15,13] = 2482.6682660008
$arr[15,14] = 0.009970507098962526
$arr[15,15] = 0.009970507098962529
$ws.Range("E2:T17").Value = $arr
Write-Output "done"
